# Add two new list entries ("3. Crates And Barrels Pack Volume 1 - Free
# Version" and "4. Cratoon FX Free") after "2. Low Poly Soldiers Demo",
# plus one extra blank paragraph, while leaving the trailing blank
# paragraphs that were already in the document untouched.
#
# The third paragraph in the original document is blank; we grow two new,
# empty paragraph marks immediately *before* it (so the blank paragraph
# itself, and the two blank paragraphs that follow it, keep being reused /
# untouched -> they stay free of any run). We then fill the two brand new
# paragraphs with the required text.

$d = $word.ActiveDocument

$target = $d.Paragraphs.Item(3)

# Insert a new empty paragraph mark right before the (still blank) 3rd
# paragraph, then give that new paragraph the first line of text.
$target.Range.InsertParagraphBefore()
$p1 = $d.Paragraphs.Item(3)
$p1.Range.Text = "3. Crates And Barrels Pack Volume 1 - Free Version"

# The original 3rd paragraph (still blank) is now the 4th paragraph.
# Insert another new empty paragraph mark right before it, and give that
# new paragraph the second line of text.
$target2 = $d.Paragraphs.Item(4)
$target2.Range.InsertParagraphBefore()
$p2 = $d.Paragraphs.Item(4)
$p2.Range.Text = "4. Cratoon FX Free"

# What remains is: ... ,"4. Cratoon FX Free", <blank>, <blank>, <blank>
# where the three blanks are the paragraph that originally held index 3
# (now pushed down to index 5) plus the two paragraphs that always
# followed it - none of them were ever split, so none of them carry a
# stray run.
